$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Row 19: columns B:K and N were blank placeholders; the new event data fills
# them in with the literal "nan" marker used throughout this sheet.
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","N")) {
    $ws.Range($col + "19").Value = "nan"
}

# Row 20: brand-new service event row appended below row 19.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "18"
$ws.Range("L20").Value = "23\10\2025"
$ws.Range("M20").Value = "883 t"
$ws.Range("O20").Value = "تم تغيير الجرئد الخلفيه (1_5_8) ومعارته"
$ws.Range("P20").Value = "الخبير"
